$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- Update row heights ---
$ws.Rows.Item(4).RowHeight = 123
$ws.Rows.Item(7).RowHeight = 98

# --- Update column widths (B=23, C=23.6640625, D=23) ---
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668
$ws.Columns.Item(3).ColumnWidth = 22.830729166666668
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668

# --- Update cell contents for rows 6-8 (new shared strings appended in this order) ---
$ws.Range("B6").Value = "Показати швидікість та зручність оформлення покупки квитків та зручність пізнання інформації"
$ws.Range("C6").Value = "Показати легку оплату за допомогою P24, MonoDirect, LiqPay, pay.Fondy.ua, GPay та Apple Pay"
$ws.Range("B7").Value = "Обирає виставу яка сподобалась, обирає дату та оплачує квиток"
$ws.Range("B8").Value = "Обирає виставу яка сподобалась, переглядає опис події, потім обирає дату та оплачує квиток"
$ws.Range("C7").Value = "Дуже уважно читає інструкцію з оплати. Перший раз вводить номер карти вручну, що не дуже безпечно, та не дуже швидко"
$ws.Range("C8").Value = "Користується підключеним екваерінгом до сервісу,ввести тільки номер телефону та підтвердити запит на оплату в мобільному банкінгу"
$ws.Range("D6").Value = "Показати календар з детальним та зручним графіком вистав, та зручним переходом до сторінки з виставою"
$ws.Range("D7").Value = "Обрала виставу та змогла переглянути в календарі, коли ще будуть постановки цієї вистави"
$ws.Range("D8").Value = "Корегує свої плани на тиждень, коли зʼявляється цікава постанова. Може обрати зручний для себе день"

# --- Update the view: scroll so row 5 is at top, select D9 ---
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("D9").Select() | Out-Null
